# Apply the diff: add a second worksheet "Table_2" with adequacy-ratio
# data, and remove the now-empty inlineStr cells (B2, A3, B37) on
# "Table_1" that the original file carried as placeholders.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Table_1: drop the empty placeholder cells -------------------------
$ws1.Range("B2").ClearContents()
$ws1.Range("A3").ClearContents()
$ws1.Range("B37").ClearContents()

# --- add the new sheet, right after Table_1 -----------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Table_2"

# Copy the header style (bold, centered, thin-bordered) already used by
# Table_1!A1 so we don't introduce a brand new style entry.
$ws1.Range("A1").Copy()
$ws2.Range("A1:D1").PasteSpecial(-4122)

$ws2.Range("A1").Value = "Əmsal"
$ws2.Range("B1").Value = "Norma (Sistem əhəmiyyətli)"
$ws2.Range("C1").Value = "Norma (Banklar istisna)"
$ws2.Range("D1").Value = "Fakt"

# Cells holding bare "NN.N%"-shaped text: Excel's autoconvert would turn a
# plain .Value assignment into a percentage number, so mark them as Text
# first and assign the literal string afterwards.
$percentCells = "B2", "C2", "D2", "B3", "C3", "D3", "D4"
foreach ($addr in $percentCells) {
    $ws2.Range($addr).NumberFormat = "@"
}

$ws2.Range("A2").Value = "9.  I dərəcəli  kapitalın  adekvatlıq əmsalı"
$ws2.Range("B2").Value = "6.0%"
$ws2.Range("C2").Value = "5.0%"
$ws2.Range("D2").Value = "12.0%"

$ws2.Range("A3").Value = "10. məcmu kapitalın  adekvatlıq  əmsalı"
$ws2.Range("B3").Value = "11.0%"
$ws2.Range("C3").Value = "9.0%"
$ws2.Range("D3").Value = "17.0%"

$ws2.Range("A4").Value = "11. Leverec əmsalı"
$ws2.Range("B4").Value = "minimum 5%"
$ws2.Range("C4").Value = "minimum 4%"
$ws2.Range("D4").Value = "7.0%"
